# Added server, write incoming data to database
# -> New "sql_type" column inserted into the register map worksheet
#    (before the existing "field" column), classifying each register's
#    output SQL storage type (INTEGER vs REAL) for the new DB-writing
#    server code.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at D (pushes field/description/bit-columns right by one).
$ws.Columns("D:D").Insert()

# Header
$ws.Range("D1").Value = "sql_type"

# Per-register sql_type values
$ws.Range("D2").Value  = "INTEGER"   # addr
$ws.Range("D3").Value  = "INTEGER"   # retries
$ws.Range("D4").Value  = "REAL"      # vbat
$ws.Range("D5").Value  = "REAL"      # temperature_internal
$ws.Range("D6").Value  = "REAL"      # si7021_humidity
$ws.Range("D7").Value  = "REAL"      # si7021_temperature
$ws.Range("D8").Value  = "REAL"      # bmp180_temperature
$ws.Range("D9").Value  = "REAL"      # bmp180_pressure
$ws.Range("D10").Value = "INTEGER"   # reed

# Match the new column widths from the authored change: sql_type (D) and
# field (E) both 18.88 wide; description (F) keeps its 71.33 width.
$ws.Columns("D:E").ColumnWidth = 18.88

# Reflect the author's final viewport/selection in the saved sheet view.
$ws.Range("D11").Select()
$ws.Application.ActiveWindow.ScrollColumn = 2
